# "smart merge on coordinates and readsheets"
# The worksheet holding the ERP coordinate/resistivity readings is renamed
# from the generic default "Sheet1" to the descriptive "testsafe".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "testsafe"
